{"js": "// Remove the three placeholder heading paragraphs (\"Heading1\", \"Heading2\",\n// \"Heading3\") that precede the trailing bookmark paragraph, leaving a single\n// empty (Normal-style) paragraph that still carries the Word-managed\n// \"_GoBack\" bookmark.\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document starts out as:\n//   [0] \"Heading1\"  (pStyle Heading1)\n//   [1] \"Heading2\"  (pStyle Heading2)\n//   [2] \"Heading3\"  (pStyle Heading3) + bookmarkStart/bookmarkEnd \"_GoBack\"\n// Delete the first two heading paragraphs outright.\nparagraphs.items[0].delete();\nparagraphs.items[1].delete();\nawait context.sync();\n\n// Remove the remaining \"Heading3\" run text without touching the paragraph\n// mark or the bookmark that follows it.\nconst searchResults = body.search(\"Heading3\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\nsearchResults.items[0].delete();\nawait context.sync();\n\n// Reset the now-empty paragraph back to the default \"Normal\" style so no\n// heading formatting (pStyle) lingers on it.\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].style = \"Normal\";\nawait context.sync();\n", "ps1": "# Remove the three placeholder heading paragraphs (\"Heading1\", \"Heading2\",\n# \"Heading3\") that precede the trailing bookmark paragraph, leaving a single\n# empty (Normal-style) paragraph that still carries the Word-managed\n# \"_GoBack\" bookmark.\n$d = $word.ActiveDocument\n\n# The document starts out as:\n#   1: \"Heading1\"  (pStyle Heading1)\n#   2: \"Heading2\"  (pStyle Heading2)\n#   3: \"Heading3\"  (pStyle Heading3) + bookmarkStart/bookmarkEnd \"_GoBack\"\n# Delete the first two heading paragraphs outright (deleting paragraph 1\n# twice shifts paragraph 3 down to index 1 each time).\n$d.Paragraphs(1).Range.Delete()\n$d.Paragraphs(1).Range.Delete()\n\n# Remove the remaining \"Heading3\" run text without touching the paragraph\n# mark or the bookmark that follows it.\n$find = $d.Content.Find\n$find.Text = \"Heading3\"\n$find.Execute() | Out-Null\nif ($find.Found) {\n    $find.Parent.Delete()\n}\n\n# Reset the now-empty paragraph back to the default \"Normal\" style so no\n# heading formatting (pStyle) lingers on it.\n$d.Paragraphs(1).Style = \"Normal\"\n"}
